# Auto-generated edit script: reorders/reshuffles paragraph content per the target diff.
$d = $word.ActiveDocument

$missing = @()

# Paragraph 6 (originally run 1)
$rng = $d.Paragraphs.Item(6).Range.Duplicate
$found = $rng.Find.Execute("Orientar os estudantes no início de sua trajetória universitária no curso de graduação em Engenharia XX na EEL-USP de modo que o estudante seja capaz de a) identificar as oportunidades acadêmicas e as particularidades do seu curso; b) reconhecer, sob acompanhamento de um tutor, eventuais dificuldades ao longo do curso e compreender mecanismos para que estas sejam superadas, conduzindo o curso com o sucesso desejado; c) desenvolver habilidades técnicas e emocionais, ampliando as perspectivas de formação profissional por meio de atividades e encontros sistematizados.", $true, $false, $false, $false, $false, $true, 1, $false, "Os cursos de engenharia física, respectivos projetos pedagógicos e seus componentes curriculares, incluindo TCC, estágio obrigatório, Projetos de Extensão Curricularizados, Atividades Acadêmicas Complementares e Atividades extracurriculares. Identificação e aderência do estudante com o curso e com a profissão escolhida. O curso superior, a transição adolescente/jovem adulto e os desafios nos projetos de vida do estudante no início da graduação. Relação entre as disciplinas e o conhecimento a ser aplicado. Competências e habilidades desenvolvidas no seu curso de engenharia. Dimensões acadêmicas, socioculturais e científicas. Diversidade e inclusão. Organização dos estudos.", 2)
if (-not $found) { $missing += "para 6 run 1" }

# Paragraph 7 (originally run 1)
$rng = $d.Paragraphs.Item(7).Range.Duplicate
$found = $rng.Find.Execute("To guide students at the beginning of their university career in the undergraduate course in Engineering XX at EEL-USP so that the student is able to: a) identify the academic opportunities and particularities of their course; b) recognize, under the supervision of a tutor, any difficulties throughout the course and understand mechanisms for them to be overcome, conducting the course with the desired success; c) develop technical and emotional skills, broadening the perspectives of professional training through systematized activities and meetings.", $true, $false, $false, $false, $false, $true, 1, $false, "Engineering physics courses, respective pedagogical projects and their curricular components, including TCC, mandatory internship, Curricular Extension Projects, Complementary Academic Activities and Extracurricular Activities. Identification and adherence of the student with the course and with the chosen profession. The college course, the adolescent/young adult transition and the challenges in the student's life projects at the beginning of the undergraduate program. Relationship between the disciplines and the knowledge to be applied. Competencies and skills developed in your engineering course. Academic, sociocultural and scientific dimensions. Diversity and inclusion. Organization of studies.", 2)
if (-not $found) { $missing += "para 7 run 1" }

# Paragraph 9 (originally run 1)
$rng = $d.Paragraphs.Item(9).Range.Duplicate
$found = $rng.Find.Execute("5817692 - Katia Cristiane Gandolpho Candioto", $true, $false, $false, $false, $false, $true, 1, $false, "Orientar os estudantes no início de sua trajetória universitária no curso de graduação em Engenharia XX na EEL-USP de modo que o estudante seja capaz de a) identificar as oportunidades acadêmicas e as particularidades do seu curso; b) reconhecer, sob acompanhamento de um tutor, eventuais dificuldades ao longo do curso e compreender mecanismos para que estas sejam superadas, conduzindo o curso com o sucesso desejado; c) desenvolver habilidades técnicas e emocionais, ampliando as perspectivas de formação profissional por meio de atividades e encontros sistematizados.", 2)
if (-not $found) { $missing += "para 9 run 1" }

# Paragraph 9 (originally run 2)
$rng = $d.Paragraphs.Item(9).Range.Duplicate
$found = $rng.Find.Execute("1176388 - Luiz Tadeu Fernandes Eleno", $true, $false, $false, $false, $false, $true, 1, $false, "Apresentação dos programas e serviços oferecidos pela USP voltados aos estudantes e das oportunidades de realizar trabalhos extracurriculares. A dinâmica das aulas, ferramentas de interação. Desenvolvimento de atividades de grupo, com objetivo de desenvolver habilidades sócio-comportamentais através de colaboração em temas do curso relacionados à profissão escolhida. Áreas de atuação do curso de engenharia, competências e habilidades a serem desenvolvidas. Interdisciplinaridade e a relação entre as disciplinas e o conhecimento a ser aplicado. Planejamento de estudos. Formas de estudar e aprender.", 2)
if (-not $found) { $missing += "para 9 run 2" }

# Paragraph 11 (originally run 1)
$rng = $d.Paragraphs.Item(11).Range.Duplicate
$found = $rng.Find.Execute("Os cursos de engenharia física, respectivos projetos pedagógicos e seus componentes curriculares, incluindo TCC, estágio obrigatório, Projetos de Extensão Curricularizados, Atividades Acadêmicas Complementares e Atividades extracurriculares. Identificação e aderência do estudante com o curso e com a profissão escolhida. O curso superior, a transição adolescente/jovem adulto e os desafios nos projetos de vida do estudante no início da graduação. Relação entre as disciplinas e o conhecimento a ser aplicado. Competências e habilidades desenvolvidas no seu curso de engenharia. Dimensões acadêmicas, socioculturais e científicas. Diversidade e inclusão. Organização dos estudos.", $true, $false, $false, $false, $false, $true, 1, $false, "Atividades realizadas na forma de dinâmicas de grupos, utilização de vídeos, textos, roda de discussão e/ou elaboração de painéis. Participação em encontros de orientação promovidos pelo Programa de Tutoria Acadêmica e a realização de atividades propostas pelo tutor/monitor/mentor, incluindo trabalhos em equipe e estudos dirigidos.", 2)
if (-not $found) { $missing += "para 11 run 1" }

# Paragraph 12 (originally run 1)
$rng = $d.Paragraphs.Item(12).Range.Duplicate
$found = $rng.Find.Execute("Engineering physics courses, respective pedagogical projects and their curricular components, including TCC, mandatory internship, Curricular Extension Projects, Complementary Academic Activities and Extracurricular Activities. Identification and adherence of the student with the course and with the chosen profession. The college course, the adolescent/young adult transition and the challenges in the student's life projects at the beginning of the undergraduate program. Relationship between the disciplines and the knowledge to be applied. Competencies and skills developed in your engineering course. Academic, sociocultural and scientific dimensions. Diversity and inclusion. Organization of studies.", $true, $false, $false, $false, $false, $true, 1, $false, "To guide students at the beginning of their university career in the undergraduate course in Engineering XX at EEL-USP so that the student is able to: a) identify the academic opportunities and particularities of their course; b) recognize, under the supervision of a tutor, any difficulties throughout the course and understand mechanisms for them to be overcome, conducting the course with the desired success; c) develop technical and emotional skills, broadening the perspectives of professional training through systematized activities and meetings.", 2)
if (-not $found) { $missing += "para 12 run 1" }

# Paragraph 14 (originally run 1)
$rng = $d.Paragraphs.Item(14).Range.Duplicate
$found = $rng.Find.Execute("Apresentação dos programas e serviços oferecidos pela USP voltados aos estudantes e das oportunidades de realizar trabalhos extracurriculares. A dinâmica das aulas, ferramentas de interação. Desenvolvimento de atividades de grupo, com objetivo de desenvolver habilidades sócio-comportamentais através de colaboração em temas do curso relacionados à profissão escolhida. Áreas de atuação do curso de engenharia, competências e habilidades a serem desenvolvidas. Interdisciplinaridade e a relação entre as disciplinas e o conhecimento a ser aplicado. Planejamento de estudos. Formas de estudar e aprender.", $true, $false, $false, $false, $false, $true, 1, $false, "Participação ativa nos encontros, apresentação de estudos/pesquisa e de trabalhos realizados durante a disciplina, colaboração e engajamento nas atividades da disciplina. O estudante deverá entregar um relatório final para a disciplina. A nota final é dada pela média ponderada das notas obtidas nas diversas atividades propostas.", 2)
if (-not $found) { $missing += "para 14 run 1" }

# Paragraph 16 (originally run 6)
$rng = $d.Paragraphs.Item(16).Range.Duplicate
$found = $rng.Find.Execute("Não se aplica..", $true, $false, $false, $false, $false, $true, 1, $false, "5817692 - Katia Cristiane Gandolpho Candioto", 2)
if (-not $found) { $missing += "para 16 run 6" }

# Paragraph 16 (originally run 4)
$rng = $d.Paragraphs.Item(16).Range.Duplicate
$found = $rng.Find.Execute("Participação ativa nos encontros, apresentação de estudos/pesquisa e de trabalhos realizados durante a disciplina, colaboração e engajamento nas atividades da disciplina. O estudante deverá entregar um relatório final para a disciplina. A nota final é dada pela média ponderada das notas obtidas nas diversas atividades propostas.", $true, $false, $false, $false, $false, $true, 1, $false, "[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001.^l[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27.^l[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.^l[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.^l[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019.", 2)
if (-not $found) { $missing += "para 16 run 4" }

# Paragraph 16 (originally run 2)
$rng = $d.Paragraphs.Item(16).Range.Duplicate
$found = $rng.Find.Execute("Atividades realizadas na forma de dinâmicas de grupos, utilização de vídeos, textos, roda de discussão e/ou elaboração de painéis. Participação em encontros de orientação promovidos pelo Programa de Tutoria Acadêmica e a realização de atividades propostas pelo tutor/monitor/mentor, incluindo trabalhos em equipe e estudos dirigidos.", $true, $false, $false, $false, $false, $true, 1, $false, "Não se aplica..", 2)
if (-not $found) { $missing += "para 16 run 2" }

# Paragraph 18 (originally run 1)
$rng = $d.Paragraphs.Item(18).Range.Duplicate
$found = $rng.Find.Execute("[1] Peddy, S. The art of mentoring – Lead, follow and get out of the way. Houston: Bullion Books, 2001.[2] Zachary, L. J. The Mentor’s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promoção do bem-estar em estudantes do ensino superior. In: Programa de Monitorização e Tutorado: oito anos a promover a integração e o sucesso académico no IST. Lisboa: IST Press, 2011. p. 19-27.[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.[5] Diretrizes Curriculares Nacionais para os cursos de graduação em Engenharia. Ministério da Educação. CNE/CES, 2019.", $true, $false, $false, $false, $false, $true, 1, $false, "1176388 - Luiz Tadeu Fernandes Eleno", 2)
if (-not $found) { $missing += "para 18 run 1" }

if ($missing.Count -gt 0) {
    Write-Output ("MISSING: " + ($missing -join ", "))
} else {
    Write-Output "All replacements applied successfully."
}
